# Generate Report for Handoff
# Refresh the "latest" handoff/handback timestamps shown across the three
# report sheets (Overview, zh-cn, de-de) to reflect a new run of the
# handback report generator.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column (D) ----------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7").Value  = "2016-25-13 04:25:29"
$wsOverview.Range("D10").Value = "2016-25-13 04:25:29"
$wsOverview.Range("D11").Value = "2016-25-13 04:25:29"
$wsOverview.Range("D12").Value = "2016-25-13 04:25:29"
$wsOverview.Range("D13").Value = "2016-25-13 04:25:29"
$wsOverview.Range("D14").Value = "2016-25-13 04:25:29"
$wsOverview.Range("D15").Value = "2016-25-13 04:25:29"
$wsOverview.Range("D16").Value = "2016-25-13 04:25:29"

# --- zh-cn sheet: "Latest Handoff Datetime" column (E) ----------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7").Value  = "2016-03-13 04:25:25"
$wsZhCn.Range("E10").Value = "2016-03-13 04:25:25"
$wsZhCn.Range("E11").Value = "2016-03-13 04:25:25"
$wsZhCn.Range("E12").Value = "2016-03-13 04:25:25"
$wsZhCn.Range("E13").Value = "2016-03-13 04:25:25"
$wsZhCn.Range("E14").Value = "2016-03-13 04:25:25"
$wsZhCn.Range("E15").Value = "2016-03-13 04:25:25"
$wsZhCn.Range("E16").Value = "2016-03-13 04:25:25"

# --- de-de sheet: "Latest Handoff Datetime" column (E) ----------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7").Value  = "2016-03-13 04:25:29"
$wsDeDe.Range("E10").Value = "2016-03-13 04:25:29"
$wsDeDe.Range("E11").Value = "2016-03-13 04:25:29"
$wsDeDe.Range("E12").Value = "2016-03-13 04:25:29"
$wsDeDe.Range("E13").Value = "2016-03-13 04:25:29"
$wsDeDe.Range("E14").Value = "2016-03-13 04:25:29"
$wsDeDe.Range("E15").Value = "2016-03-13 04:25:29"
$wsDeDe.Range("E16").Value = "2016-03-13 04:25:29"
